$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$TextValue)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $TextValue
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '51.075.07'
Set-TextValue $ws.Range("E2") '  -0.01%  '
Set-TextValue $ws.Range("D3") '2.959.01'
Set-TextValue $ws.Range("E3") '  +0.43%  '
Set-TextValue $ws.Range("D4") '0.999'
Set-TextValue $ws.Range("E4") '  +0.09%  '
Set-TextValue $ws.Range("D5") '380.06'
Set-TextValue $ws.Range("E5") '  +1.53%  '
Set-TextValue $ws.Range("D6") '102.18'
Set-TextValue $ws.Range("E6") '  -0.23%  '
Set-TextValue $ws.Range("D7") '0.545'
Set-TextValue $ws.Range("E7") '  +1.78%  '
Set-TextValue $ws.Range("E8") '  +0.02%  '
Set-TextValue $ws.Range("E9") '  +1.50%  '
Set-TextValue $ws.Range("D10") '36.54'
Set-TextValue $ws.Range("E10") '  +0.33%  '
Set-TextValue $ws.Range("E11") '  -1.27%  '
Set-TextValue $ws.Range("D12") '0.0854'
Set-TextValue $ws.Range("E12") '  +2.06%  '
Set-TextValue $ws.Range("B13") 'Polkadot'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D13") '7.82'
Set-TextValue $ws.Range("E13") '  +6.36%  '
Set-TextValue $ws.Range("B14") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D14") '3.426.81'
Set-TextValue $ws.Range("E14") '  +0.70%  '
Set-TextValue $ws.Range("D15") '18.34'
Set-TextValue $ws.Range("E15") '  +2.57%  '
Set-TextValue $ws.Range("D16") '11.61'
Set-TextValue $ws.Range("E16") '  +61.35%  '
Set-TextValue $ws.Range("D17") '2.948.10'
Set-TextValue $ws.Range("E17") '  +0.27%  '
Set-TextValue $ws.Range("D18") '0.998'
Set-TextValue $ws.Range("E18") '  +1.77%  '
Set-TextValue $ws.Range("D19") '51.142.43'
Set-TextValue $ws.Range("E19") '  +0.39%  '
Set-TextValue $ws.Range("D20") '3.12'
Set-TextValue $ws.Range("E20") '  -1.01%  '
Set-TextValue $ws.Range("D21") '12.39'
Set-TextValue $ws.Range("E21") '  -2.12%  '
Set-TextValue $ws.Range("E22") '  +0.36%  '
Set-TextValue $ws.Range("B23") 'PancakeSwap'
Set-TextValue $ws.Range("C23") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D23") '3.30'
Set-TextValue $ws.Range("E23") '  +14.67%  '
Set-TextValue $ws.Range("B24") 'Litecoin'
Set-TextValue $ws.Range("C24") 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D24") '70.06'
Set-TextValue $ws.Range("E24") '  +2.52%  '
Set-TextValue $ws.Range("E25") '  +0.91%  '
Set-TextValue $ws.Range("D26") '7.89'
Set-TextValue $ws.Range("E26") '  -7.55%  '
Set-TextValue $ws.Range("E27") '  -0.04%  '
Set-TextValue $ws.Range("D28") '7.16'
Set-TextValue $ws.Range("E28") '  -11.36%  '
Set-TextValue $ws.Range("D29") '0.165'
Set-TextValue $ws.Range("E29") '  -1.66%  '
Set-TextValue $ws.Range("D30") '25.85'
Set-TextValue $ws.Range("E30") '  +0.80%  '
Set-TextValue $ws.Range("E31") '  -3.21%  '
Set-TextValue $ws.Range("D32") '10.31'
Set-TextValue $ws.Range("E32") '  +4.36%  '
Set-TextValue $ws.Range("B33") 'InjectiveProtocol'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D33") '34.38'
Set-TextValue $ws.Range("E33") '  +2.02%  '
Set-TextValue $ws.Range("B34") 'OKB'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D34") '51.12'
Set-TextValue $ws.Range("E34") '  +0.49%  '
Set-TextValue $ws.Range("E35") '  +1.94%  '
Set-TextValue $ws.Range("E36") '  -2.74%  '
Set-TextValue $ws.Range("E37") '  +0.13%  '
Set-TextValue $ws.Range("D38") '3.25'
Set-TextValue $ws.Range("E38") '  +8.65%  '
Set-TextValue $ws.Range("D39") '0.116'
Set-TextValue $ws.Range("E39") '  +0.60%  '
Set-TextValue $ws.Range("B40") 'Celestia'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D40") '16.55'
Set-TextValue $ws.Range("E40") '  +0.62%  '
Set-TextValue $ws.Range("B41") 'ARBITRUM'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D41") '1.83'
Set-TextValue $ws.Range("E41") '  +2.08%  '
Set-TextValue $ws.Range("D42") '2.50'
Set-TextValue $ws.Range("E42") '  -1.92%  '
Set-TextValue $ws.Range("D43") '124.71'
Set-TextValue $ws.Range("E43") '  +3.58%  '
Set-TextValue $ws.Range("D44") '21.55'
Set-TextValue $ws.Range("E44") '  +2.41%  '
Set-TextValue $ws.Range("E45") '  +8.17%  '
Set-TextValue $ws.Range("B46") 'WEMIXToken'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D46") '2.02'
Set-TextValue $ws.Range("E46") '  -0.91%  '
Set-TextValue $ws.Range("D47") '2.37'
Set-TextValue $ws.Range("E47") '  +3.06%  '
Set-TextValue $ws.Range("B48") 'TheGraph'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D48") '0.271'
Set-TextValue $ws.Range("E48") '  -5.13%  '
Set-TextValue $ws.Range("D49") '2.048.57'
Set-TextValue $ws.Range("E49") '  +3.72%  '
Set-TextValue $ws.Range("D50") '0.0321'
Set-TextValue $ws.Range("E50") '  -6.46%  '
Set-TextValue $ws.Range("E51") '  +7.15%  '
